# TaxPayerSSNNoMatch.xlsx - RAD Phase 3 test data refresh for "Estate Tax".
# A later Katalon run (Fri Feb 02 19:43-19:44 EST 2024) overwrote the three
# existing "Personal Income Tax" timestamps and appended two brand-new rows
# (6 and 7) recording a Pass/Y result for the "Estate Tax" tax type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh the Date column for the existing Personal Income Tax rows ---
$ws.Range("B2").Value = "Fri Feb 02 19:43:39 EST 2024"
$ws.Range("B3").Value = "Fri Feb 02 19:43:54 EST 2024"
$ws.Range("B5").Value = "Fri Feb 02 19:44:07 EST 2024"

# --- Row 6: new Estate Tax result - "Existing Liability w/Notice Number" ---
# Columns D6/E6 already held "Existing Liability w/Notice Number" / "Estate Tax";
# fill in the now-populated Result/Date/Execute columns to match rows 2-5.
$ws.Range("A6").Value = "Pass"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = "Fri Feb 02 19:44:19 EST 2024"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "Y"

# --- Row 7: new Estate Tax result - "New Tax Return Amount Due" ---
$ws.Range("A7").Value = "Pass"
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Value = "Fri Feb 02 19:44:33 EST 2024"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "Y"

# Selection left on C7, matching the author's final cursor position.
$ws.Range("C7").Select()
